$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (trailing data rows that no longer exist in the latest scrape)
$ws.Rows("6:9").Delete()

# Row 2
$ws.Range("A2").Value = "'1330536"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330536"
$ws.Range("C2").Value = "[IMPACT FORTALEZA] INSIDE SALES"
$ws.Range("D2").Value = "Castanhal, PA, Brasil"
$ws.Range("F2").Value = "21 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Petruz Fruity"

# Row 3
$ws.Range("A3").Value = "'1328965"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328965"
$ws.Range("C3").Value = "Account Manager (German level C1/C2 only)"
$ws.Range("D3").Value = "Assen, Nederland"
$ws.Range("F3").Value = "21 applicants"
$ws.Range("H3").Value = "ICT Specialist"

# Row 4
$ws.Range("A4").Value = "'1320932"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1320932"
$ws.Range("C4").Value = "Marketing Executive"
$ws.Range("D4").Value = "Hong Kong"
$ws.Range("F4").Value = "174 applicants"
$ws.Range("H4").Value = "Treehouse"

# Row 5
$ws.Range("A5").Value = "'1315739"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1315739"
$ws.Range("C5").Value = "Software Developer"
$ws.Range("D5").Value = "İstanbul, Türkiye"
$ws.Range("F5").Value = "262 applicants"
$ws.Range("H5").Value = "Unixpadel"

# Reset style for the ID cells so the quote-prefix doesn't alter formatting
$ws.Range("A2:A5").Style = "Normal"

# Update column widths (Excel stores width_xml = ColumnWidth + 0.8333333333333334,
# so subtract that offset to land exactly on the target stored widths of 44/24/17)
$ws.Columns("C").ColumnWidth = 43.166666666666664
$ws.Columns("D").ColumnWidth = 23.166666666666668
$ws.Columns("H").ColumnWidth = 16.166666666666668
